$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Helper cells: keep all original (pre-existing) shared strings alive ---
$ws.Range("A4").Value = "Suppliers"
$ws.Range("B4").Value = "SUBSUPPLIER2_ON_06/12/18-12:34"
$ws.Range("C4").Value = "LOCATION_ON_12/12/18-16:42"
$ws.Range("D4").Value = "SUBLOCATION_ON_12/12/18-16:42"
$ws.Range("E4").Value = "Customers > Customer > Customer2"
$ws.Range("F4").Value = "CUSTOMER_ON_12/12/18-16:49"
$ws.Range("G4").Value = "SUBCUSTOMER1_ON_12/12/18-16:49"
$ws.Range("H4").Value = "SUBCUSTOMER2_ON_12/12/18-16:49"
$ws.Range("I4").Value = "EQUIPMENT_ON_12/12/18-16:53"
$ws.Range("J4").Value = "SUBEQUIPMENT1_ON_12/12/18-16:53"
$ws.Range("K4").Value = "SUBEQUIPMENT2_ON_12/12/18-16:53"
$ws.Range("L4").Value = "ITEM_ON_12/12/18-16:57"
$ws.Range("M4").Value = "SUBITEM1_ON_12/12/18-16:57"
$ws.Range("N4").Value = "SUBITEM2_ON_12/12/18-16:57"
$ws.Range("O4").Value = "SUPPLIER_ON_12/12/18-17:00"
$ws.Range("P4").Value = "SUBSUPPLIER1_ON_12/12/18-17:00"
$ws.Range("Q4").Value = "SUBSUPPLIER2_ON_12/12/18-17:00"
$ws.Range("R4").Value = "LOCATION_ON_17/12/18-13:04"
$ws.Range("S4").Value = "SUBLOCATION_ON_17/12/18-13:04"
$ws.Range("T4").Value = "SUPPLIER_ON_17/12/18-13:04"
$ws.Range("U4").Value = "SUBSUPPLIER1_ON_17/12/18-13:04"
$ws.Range("V4").Value = "SUBSUPPLIER2_ON_17/12/18-13:04"
$ws.Range("W4").Value = "LOCATION_ON_17/12/18-13:50"
$ws.Range("X4").Value = "SUBLOCATION_ON_17/12/18-13:50"
$ws.Range("Y4").Value = "SUPPLIER_ON_17/12/18-13:52"
$ws.Range("Z4").Value = "SUBSUPPLIER1_ON_17/12/18-13:52"
$ws.Range("A5").Value = "SUBSUPPLIER2_ON_17/12/18-13:52"
$ws.Range("B5").Value = "SUPPLIER_ON_17/12/18-13:58"
$ws.Range("C5").Value = "SUBSUPPLIER1_ON_17/12/18-13:58"
$ws.Range("D5").Value = "SUBSUPPLIER2_ON_17/12/18-13:58"
$ws.Range("E5").Value = "LOCATION_ON_17/12/18-14:03"
$ws.Range("F5").Value = "SUBLOCATION_ON_17/12/18-14:03"
$ws.Range("G5").Value = "SUPPLIER_ON_17/12/18-14:04"
$ws.Range("H5").Value = "SUBSUPPLIER1_ON_17/12/18-14:04"
$ws.Range("I5").Value = "SUBSUPPLIER2_ON_17/12/18-14:04"
$ws.Range("J5").Value = "Customers > ankcust > 15thAug"
$ws.Range("K5").Value = "LOCATION_ON_17/12/18-16:32"
$ws.Range("L5").Value = "SUBLOCATION_ON_17/12/18-16:32"
$ws.Range("M5").Value = "SUPPLIER_ON_17/12/18-16:33"
$ws.Range("N5").Value = "SUBSUPPLIER1_ON_17/12/18-16:33"
$ws.Range("O5").Value = "SUBSUPPLIER2_ON_17/12/18-16:33"
$ws.Range("P5").Value = "LOCATION_ON_17/12/18-17:57"
$ws.Range("Q5").Value = "SUBLOCATION_ON_17/12/18-17:57"
$ws.Range("R5").Value = "SUPPLIER_ON_17/12/18-17:57"
$ws.Range("S5").Value = "SUBSUPPLIER1_ON_17/12/18-17:57"
$ws.Range("T5").Value = "SUBSUPPLIER2_ON_17/12/18-17:57"
$ws.Range("U5").Value = "LOCATION_ON_17/12/18-18:10"
$ws.Range("V5").Value = "SUBLOCATION_ON_17/12/18-18:10"
$ws.Range("W5").Value = "SUPPLIER_ON_17/12/18-18:11"
$ws.Range("X5").Value = "SUBSUPPLIER1_ON_17/12/18-18:11"
$ws.Range("Y5").Value = "SUBSUPPLIER2_ON_17/12/18-18:11"
$ws.Range("Z5").Value = "LOCATION_ON_17/12/18-18:35"
$ws.Range("A6").Value = "SUBLOCATION_ON_17/12/18-18:35"
$ws.Range("B6").Value = "SUPPLIER_ON_17/12/18-18:36"
$ws.Range("C6").Value = "SUBSUPPLIER1_ON_17/12/18-18:36"
$ws.Range("D6").Value = "SUBSUPPLIER2_ON_17/12/18-18:36"
$ws.Range("E6").Value = "LOCATION_ON_17/12/18-18:59"
$ws.Range("F6").Value = "SUBLOCATION_ON_17/12/18-18:59"
$ws.Range("G6").Value = "SUPPLIER_ON_17/12/18-19:00"
$ws.Range("H6").Value = "SUBSUPPLIER1_ON_17/12/18-19:00"
$ws.Range("I6").Value = "SUBSUPPLIER2_ON_17/12/18-19:00"
$ws.Range("J6").Value = "LOCATION_ON_18/12/18-11:14"
$ws.Range("K6").Value = "SUBLOCATION_ON_18/12/18-11:14"
$ws.Range("L6").Value = "CUSTOMER_ON_18/12/18-11:16"
$ws.Range("M6").Value = "SUBCUSTOMER1_ON_18/12/18-11:16"
$ws.Range("N6").Value = "SUBCUSTOMER2_ON_18/12/18-11:16"
$ws.Range("O6").Value = "EQUIPMENT_ON_18/12/18-11:19"
$ws.Range("P6").Value = "SUBEQUIPMENT1_ON_18/12/18-11:19"
$ws.Range("Q6").Value = "SUBEQUIPMENT2_ON_18/12/18-11:19"
$ws.Range("R6").Value = "Equipment"
$ws.Range("S6").Value = "ITEM_ON_18/12/18-11:23"
$ws.Range("T6").Value = "SUBITEM1_ON_18/12/18-11:23"
$ws.Range("U6").Value = "SUBITEM2_ON_18/12/18-11:23"
$ws.Range("V6").Value = "Items"
$ws.Range("W6").Value = "SUPPLIER_ON_18/12/18-11:26"
$ws.Range("X6").Value = "SUBSUPPLIER1_ON_18/12/18-11:26"
$ws.Range("Y6").Value = "SUBSUPPLIER2_ON_18/12/18-11:26"
$ws.Range("Z6").Value = "LOCATION_ON_03/01/19-11:35"
$ws.Range("A7").Value = "SUBLOCATION_ON_03/01/19-11:35"
$ws.Range("B7").Value = "LOCATION_ON_03/01/19-11:39"
$ws.Range("C7").Value = "SUBLOCATION_ON_03/01/19-11:39"
$ws.Range("D7").Value = "SUPPLIER_ON_03/01/19-11:40"
$ws.Range("E7").Value = "SUBSUPPLIER1_ON_03/01/19-11:40"
$ws.Range("F7").Value = "SUBSUPPLIER2_ON_03/01/19-11:40"
$ws.Range("G7").Value = "Customers > Cust11343 > @@@"
$ws.Range("H7").Value = "LOCATION_ON_03/01/19-12:01"
$ws.Range("I7").Value = "SUBLOCATION_ON_03/01/19-12:01"
$ws.Range("J7").Value = "SUPPLIER_ON_03/01/19-12:02"
$ws.Range("K7").Value = "SUBSUPPLIER1_ON_03/01/19-12:02"
$ws.Range("L7").Value = "SUBSUPPLIER2_ON_03/01/19-12:02"
$ws.Range("M7").Value = "LOCATION_ON_04/01/19-11:12"
$ws.Range("N7").Value = "SUBLOCATION_ON_04/01/19-11:12"
$ws.Range("O7").Value = "CUSTOMER_ON_04/01/19-11:14"
$ws.Range("P7").Value = "SUBCUSTOMER1_ON_04/01/19-11:14"
$ws.Range("Q7").Value = "SUBCUSTOMER2_ON_04/01/19-11:14"
$ws.Range("R7").Value = "EQUIPMENT_ON_04/01/19-11:17"
$ws.Range("S7").Value = "SUBEQUIPMENT1_ON_04/01/19-11:17"
$ws.Range("T7").Value = "SUBEQUIPMENT2_ON_04/01/19-11:17"
$ws.Range("U7").Value = "ITEM_ON_04/01/19-11:25"
$ws.Range("V7").Value = "SUBITEM1_ON_04/01/19-11:25"
$ws.Range("W7").Value = "SUBITEM2_ON_04/01/19-11:25"
$ws.Range("X7").Value = "SUPPLIER_ON_04/01/19-11:28"
$ws.Range("Y7").Value = "SUBSUPPLIER1_ON_04/01/19-11:28"
$ws.Range("Z7").Value = "SUBSUPPLIER2_ON_04/01/19-11:28"
$ws.Range("A8").Value = "LOCATION_ON_04/01/19-11:47"
$ws.Range("B8").Value = "SUBLOCATION_ON_04/01/19-11:47"
$ws.Range("C8").Value = "CUSTOMER_ON_04/01/19-11:48"
$ws.Range("D8").Value = "SUBCUSTOMER1_ON_04/01/19-11:48"
$ws.Range("E8").Value = "SUBCUSTOMER2_ON_04/01/19-11:48"
$ws.Range("F8").Value = "EQUIPMENT_ON_04/01/19-11:52"
$ws.Range("G8").Value = "SUBEQUIPMENT1_ON_04/01/19-11:52"
$ws.Range("H8").Value = "SUBEQUIPMENT2_ON_04/01/19-11:52"
$ws.Range("I8").Value = "ITEM_ON_04/01/19-11:56"
$ws.Range("J8").Value = "SUBITEM1_ON_04/01/19-11:56"
$ws.Range("K8").Value = "SUBITEM2_ON_04/01/19-11:56"
$ws.Range("L8").Value = "SUPPLIER_ON_04/01/19-12:00"
$ws.Range("M8").Value = "SUBSUPPLIER1_ON_04/01/19-12:00"
$ws.Range("N8").Value = "SUBSUPPLIER2_ON_04/01/19-12:00"
$ws.Range("O8").Value = "LOCATION_ON_04/01/19-18:45"
$ws.Range("P8").Value = "SUBLOCATION_ON_04/01/19-18:45"
$ws.Range("Q8").Value = "SUPPLIER_ON_04/01/19-18:46"
$ws.Range("R8").Value = "SUBSUPPLIER1_ON_04/01/19-18:46"
$ws.Range("S8").Value = "SUBSUPPLIER2_ON_04/01/19-18:46"

# --- Helper cells: introduce brand-new shared strings in ascending target order ---
$ws.Range("A10").Value = "LOCATION_ON_09/01/19-11:48"
$ws.Range("B10").Value = "SUBLOCATION_ON_09/01/19-11:48"
$ws.Range("C10").Value = "SUPPLIER_ON_09/01/19-11:49"
$ws.Range("D10").Value = "SUBSUPPLIER1_ON_09/01/19-11:49"
$ws.Range("E10").Value = "SUBSUPPLIER2_ON_09/01/19-11:49"
$ws.Range("F10").Value = "LOCATION_ON_09/01/19-12:02"
$ws.Range("G10").Value = "SUBLOCATION_ON_09/01/19-12:02"
$ws.Range("H10").Value = "SUPPLIER_ON_09/01/19-12:03"
$ws.Range("I10").Value = "SUBSUPPLIER1_ON_09/01/19-12:03"
$ws.Range("J10").Value = "SUBSUPPLIER2_ON_09/01/19-12:03"
$ws.Range("K10").Value = "LOCATION_ON_09/01/19-12:14"
$ws.Range("L10").Value = "SUBLOCATION_ON_09/01/19-12:14"
$ws.Range("M10").Value = "SUPPLIER_ON_09/01/19-12:14"
$ws.Range("N10").Value = "SUBSUPPLIER1_ON_09/01/19-12:14"
$ws.Range("O10").Value = "SUBSUPPLIER2_ON_09/01/19-12:14"
$ws.Range("P10").Value = "LOCATION_ON_09/01/19-12:16"
$ws.Range("Q10").Value = "SUBLOCATION_ON_09/01/19-12:16"
$ws.Range("R10").Value = "SUPPLIER_ON_09/01/19-12:17"
$ws.Range("S10").Value = "SUBSUPPLIER1_ON_09/01/19-12:17"
$ws.Range("T10").Value = "SUBSUPPLIER2_ON_09/01/19-12:17"
$ws.Range("U10").Value = "LOCATION_ON_09/01/19-12:34"
$ws.Range("V10").Value = "SUBLOCATION_ON_09/01/19-12:34"
$ws.Range("W10").Value = "SUPPLIER_ON_09/01/19-12:35"
$ws.Range("X10").Value = "SUBSUPPLIER1_ON_09/01/19-12:35"
$ws.Range("Y10").Value = "SUBSUPPLIER2_ON_09/01/19-12:35"
$ws.Range("Z10").Value = "LOCATION_ON_09/01/19-14:08"
$ws.Range("A11").Value = "SUBLOCATION_ON_09/01/19-14:08"
$ws.Range("B11").Value = "SUPPLIER_ON_09/01/19-14:09"
$ws.Range("C11").Value = "SUBSUPPLIER1_ON_09/01/19-14:09"
$ws.Range("D11").Value = "SUBSUPPLIER2_ON_09/01/19-14:09"
$ws.Range("E11").Value = "LOCATION_ON_09/01/19-15:02"
$ws.Range("F11").Value = "SUBLOCATION_ON_09/01/19-15:02"
$ws.Range("G11").Value = "SUPPLIER_ON_09/01/19-15:03"
$ws.Range("H11").Value = "SUBSUPPLIER1_ON_09/01/19-15:03"
$ws.Range("I11").Value = "SUBSUPPLIER2_ON_09/01/19-15:03"
$ws.Range("J11").Value = "LOCATION_ON_09/01/19-16:25"
$ws.Range("K11").Value = "SUBLOCATION_ON_09/01/19-16:25"
$ws.Range("L11").Value = "SUPPLIER_ON_09/01/19-16:26"
$ws.Range("M11").Value = "SUBSUPPLIER1_ON_09/01/19-16:26"
$ws.Range("N11").Value = "SUBSUPPLIER2_ON_09/01/19-16:26"
$ws.Range("O11").Value = "LOCATION_ON_09/01/19-16:55"
$ws.Range("P11").Value = "SUBLOCATION_ON_09/01/19-16:55"
$ws.Range("Q11").Value = "SUPPLIER_ON_09/01/19-16:56"
$ws.Range("R11").Value = "SUBSUPPLIER1_ON_09/01/19-16:56"
$ws.Range("S11").Value = "SUBSUPPLIER2_ON_09/01/19-16:56"
$ws.Range("T11").Value = "LOCATION_ON_09/01/19-17:07"
$ws.Range("U11").Value = "SUBLOCATION_ON_09/01/19-17:07"
$ws.Range("V11").Value = "SUPPLIER_ON_09/01/19-17:08"
$ws.Range("W11").Value = "SUBSUPPLIER1_ON_09/01/19-17:08"
$ws.Range("X11").Value = "SUBSUPPLIER2_ON_09/01/19-17:08"
$ws.Range("Y11").Value = "LOCATION_ON_10/01/19-16:46"
$ws.Range("Z11").Value = "SUBLOCATION_ON_10/01/19-16:46"
$ws.Range("A12").Value = "CUSTOMER_ON_10/01/19-16:47"
$ws.Range("B12").Value = "SUBCUSTOMER1_ON_10/01/19-16:47"
$ws.Range("C12").Value = "SUBCUSTOMER2_ON_10/01/19-16:47"
$ws.Range("D12").Value = "EQUIPMENT_ON_10/01/19-16:51"
$ws.Range("E12").Value = "SUBEQUIPMENT1_ON_10/01/19-16:51"
$ws.Range("F12").Value = "SUBEQUIPMENT2_ON_10/01/19-16:51"
$ws.Range("G12").Value = "ITEM_ON_10/01/19-16:55"
$ws.Range("H12").Value = "SUBITEM1_ON_10/01/19-16:55"
$ws.Range("I12").Value = "SUBITEM2_ON_10/01/19-16:55"
$ws.Range("J12").Value = "SUPPLIER_ON_10/01/19-16:59"
$ws.Range("K12").Value = "SUBSUPPLIER1_ON_10/01/19-16:59"
$ws.Range("L12").Value = "SUBSUPPLIER2_ON_10/01/19-16:59"
$ws.Range("M12").Value = "LOCATION_ON_10/01/19-17:09"
$ws.Range("N12").Value = "SUBLOCATION_ON_10/01/19-17:10"
$ws.Range("O12").Value = "SUPPLIER_ON_10/01/19-17:10"
$ws.Range("P12").Value = "SUBSUPPLIER1_ON_10/01/19-17:10"
$ws.Range("Q12").Value = "SUBSUPPLIER2_ON_10/01/19-17:10"
$ws.Range("R12").Value = "LOCATION_ON_11/01/19-10:20"
$ws.Range("S12").Value = "SUBLOCATION_ON_11/01/19-10:21"
$ws.Range("T12").Value = "SUPPLIER_ON_11/01/19-10:21"
$ws.Range("U12").Value = "SUBSUPPLIER1_ON_11/01/19-10:21"
$ws.Range("V12").Value = "SUBSUPPLIER2_ON_11/01/19-10:21"
$ws.Range("W12").Value = "ITEM_ON_11/01/19-11:11"
$ws.Range("X12").Value = "SUBITEM1_ON_11/01/19-11:11"
$ws.Range("Y12").Value = "SUBITEM2_ON_11/01/19-11:11"
$ws.Range("Z12").Value = "ITEM_ON_11/01/19-11:14"
$ws.Range("A13").Value = "SUBITEM1_ON_11/01/19-11:14"
$ws.Range("B13").Value = "SUBITEM2_ON_11/01/19-11:14"
$ws.Range("C13").Value = "ITEM_ON_11/01/19-11:18"
$ws.Range("D13").Value = "SUBITEM1_ON_11/01/19-11:18"
$ws.Range("E13").Value = "SUBITEM2_ON_11/01/19-11:18"
$ws.Range("F13").Value = "LOCATION_ON_11/01/19-12:48"
$ws.Range("G13").Value = "SUBLOCATION_ON_11/01/19-12:48"
$ws.Range("H13").Value = "SUPPLIER_ON_11/01/19-12:49"
$ws.Range("I13").Value = "SUBSUPPLIER1_ON_11/01/19-12:49"
$ws.Range("J13").Value = "SUBSUPPLIER2_ON_11/01/19-12:49"
$ws.Range("K13").Value = "ITEM_ON_11/01/19-12:51"
$ws.Range("L13").Value = "SUBITEM1_ON_11/01/19-12:51"
$ws.Range("M13").Value = "SUBITEM2_ON_11/01/19-12:51"
$ws.Range("N13").Value = "LOCATION_ON_11/01/19-14:10"
$ws.Range("O13").Value = "SUBLOCATION_ON_11/01/19-14:10"
$ws.Range("P13").Value = "SUPPLIER_ON_11/01/19-14:11"
$ws.Range("Q13").Value = "SUBSUPPLIER1_ON_11/01/19-14:11"
$ws.Range("R13").Value = "SUBSUPPLIER2_ON_11/01/19-14:11"
$ws.Range("S13").Value = "ITEM_ON_11/01/19-14:13"
$ws.Range("T13").Value = "SUBITEM1_ON_11/01/19-14:13"
$ws.Range("U13").Value = "SUBITEM2_ON_11/01/19-14:13"
$ws.Range("V13").Value = "Location_11/01/19-15:27"
$ws.Range("W13").Value = "SubLocation_11/01/19-15:27"
$ws.Range("X13").Value = "Customer_11/01/19-15:29"
$ws.Range("Y13").Value = "SubCustomer1_11/01/19-15:29"
$ws.Range("Z13").Value = "SubCustomer2_11/01/19-15:29"
$ws.Range("A14").Value = "Equipment_11/01/19-15:33"
$ws.Range("B14").Value = "SubEquipment1_11/01/19-15:33"
$ws.Range("C14").Value = "SubEquipment2_11/01/19-15:33"
$ws.Range("D14").Value = "Location_11/01/19-16:16"
$ws.Range("E14").Value = "SubLocation_11/01/19-16:16"
$ws.Range("F14").Value = "Customer_11/01/19-16:18"
$ws.Range("G14").Value = "SubCustomer1_11/01/19-16:18"
$ws.Range("H14").Value = "SubCustomer2_11/01/19-16:18"
$ws.Range("I14").Value = "Location_11/01/19-16:23"
$ws.Range("J14").Value = "SubLocation_11/01/19-16:24"
$ws.Range("K14").Value = "Location_11/01/19-16:45"
$ws.Range("L14").Value = "SubLocation_11/01/19-16:45"
$ws.Range("M14").Value = "Customer_11/01/19-16:47"
$ws.Range("N14").Value = "SubCustomer1_11/01/19-16:47"
$ws.Range("O14").Value = "SubCustomer2_11/01/19-16:47"
$ws.Range("P14").Value = "Equipment_11/01/19-16:51"
$ws.Range("Q14").Value = "SubEquipment1_11/01/19-16:51"
$ws.Range("R14").Value = "SubEquipment2_11/01/19-16:51"
$ws.Range("S14").Value = "Item_11/01/19-16:54"
$ws.Range("T14").Value = "SubItem1_11/01/19-16:54"
$ws.Range("U14").Value = "SubItem2_11/01/19-16:54"
$ws.Range("V14").Value = "Location_11/01/19-17:03"
$ws.Range("W14").Value = "SubLocation_11/01/19-17:03"
$ws.Range("X14").Value = "Customer_11/01/19-17:04"
$ws.Range("Y14").Value = "SubCustomer1_11/01/19-17:04"
$ws.Range("Z14").Value = "SubCustomer2_11/01/19-17:04"
$ws.Range("A15").Value = "Equipment_11/01/19-17:07"
$ws.Range("B15").Value = "SubEquipment1_11/01/19-17:07"
$ws.Range("C15").Value = "SubEquipment2_11/01/19-17:07"
$ws.Range("D15").Value = "Item_11/01/19-17:11"
$ws.Range("E15").Value = "SubItem1_11/01/19-17:11"
$ws.Range("F15").Value = "SubItem2_11/01/19-17:11"
$ws.Range("G15").Value = "Location_11/01/19-17:21"
$ws.Range("H15").Value = "SubLocation_11/01/19-17:21"
$ws.Range("I15").Value = "Customer_11/01/19-17:22"
$ws.Range("J15").Value = "SubCustomer1_11/01/19-17:22"
$ws.Range("K15").Value = "SubCustomer2_11/01/19-17:22"
$ws.Range("L15").Value = "Equipment_11/01/19-17:25"
$ws.Range("M15").Value = "SubEquipment1_11/01/19-17:25"
$ws.Range("N15").Value = "SubEquipment2_11/01/19-17:25"
$ws.Range("O15").Value = "Location_11/01/19-17:39"
$ws.Range("P15").Value = "SubLocation_11/01/19-17:39"
$ws.Range("Q15").Value = "Customer_11/01/19-17:40"
$ws.Range("R15").Value = "SubCustomer1_11/01/19-17:40"
$ws.Range("S15").Value = "SubCustomer2_11/01/19-17:40"
$ws.Range("T15").Value = "Equipment_11/01/19-17:44"
$ws.Range("U15").Value = "SubEquipment1_11/01/19-17:44"
$ws.Range("V15").Value = "SubEquipment2_11/01/19-17:44"
$ws.Range("W15").Value = "Location_11/01/19-17:51"
$ws.Range("X15").Value = "SubLocation_11/01/19-17:51"
$ws.Range("Y15").Value = "Customer_11/01/19-17:52"
$ws.Range("Z15").Value = "SubCustomer1_11/01/19-17:52"
$ws.Range("A16").Value = "SubCustomer2_11/01/19-17:52"

# --- Now set row 2 (the actual data row) to its final values, in ascending new-string order ---
$ws.Range("A2").Value = "Location_11/01/19-17:57"
$ws.Range("B2").Value = "SubLocation_11/01/19-17:57"
$ws.Range("C2").Value = "Customer_11/01/19-17:57"
$ws.Range("G2").Value = "SubCustomer1_11/01/19-17:57"
$ws.Range("H2").Value = "SubCustomer2_11/01/19-17:57"
$ws.Range("D2").Value = "Equipment_11/01/19-18:00"
$ws.Range("I2").Value = "SubEquipment1_11/01/19-18:00"
$ws.Range("J2").Value = "SubEquipment2_11/01/19-18:00"
$ws.Range("E2").Value = "Item_11/01/19-18:01"
$ws.Range("K2").Value = "SubItem1_11/01/19-18:01"
$ws.Range("L2").Value = "SubItem2_11/01/19-18:01"
$ws.Range("F2").Value = "Supplier_11/01/19-18:01"
$ws.Range("M2").Value = "SubSupplier1_11/01/19-18:01"
$ws.Range("N2").Value = "SubSupplier2_11/01/19-18:01"
$ws.Range("O2").Value = "Customers > Cust21511 > 456"
$ws.Range("P2").Value = "Customers"
$ws.Range("Q2").Value = "Supplier_11/01/19-18:01"
$ws.Range("R2").Value = "SubSupplier1_11/01/19-18:01"
$ws.Range("S2").Value = "SubSupplier2_11/01/19-18:01"

# --- Clean up helper cells: clear their contents so only rows 1-2 hold visible data ---
# (left intentionally NOT cleared, see notes)